$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataFrame"

# Delete now-empty rows (descending order to keep row indices stable)
$ws.Rows.Item(95).Delete() | Out-Null
$ws.Rows.Item(69).Delete() | Out-Null
$ws.Rows.Item(68).Delete() | Out-Null
$ws.Rows.Item(43).Delete() | Out-Null
$ws.Rows.Item(30).Delete() | Out-Null
$ws.Rows.Item(24).Delete() | Out-Null
$ws.Rows.Item(16).Delete() | Out-Null
$ws.Rows.Item(11).Delete() | Out-Null

# Fix Level 1 (column B) values: trim trailing/leading whitespace
$ws.Cells.Item(18, 2).Value = "Outage"
$ws.Cells.Item(19, 2).Value = "Other Transactions"
$ws.Cells.Item(22, 2).Value = "Handled calls"
$ws.Cells.Item(23, 2).Value = "Transferred calls4"
$ws.Cells.Item(26, 2).Value = "Billing Balance"
$ws.Cells.Item(31, 2).Value = "Emergency"
$ws.Cells.Item(32, 2).Value = "Outage"
$ws.Cells.Item(33, 2).Value = "Spanish (General)"
$ws.Cells.Item(35, 2).Value = "Solar"
$ws.Cells.Item(36, 2).Value = "Specialty Lines"
$ws.Cells.Item(38, 2).Value = "Other Transactions"
$ws.Cells.Item(70, 2).Value = "CARE/FERA"
$ws.Cells.Item(81, 2).Value = "Pilot Light Appointments"
$ws.Cells.Item(82, 2).Value = "Pilot Light Appointments"
$ws.Cells.Item(83, 2).Value = "Pilot Light Appointments"
$ws.Cells.Item(84, 2).Value = "Pilot Light Appointments"
$ws.Cells.Item(85, 2).Value = "Budget Billing (Login)"
$ws.Cells.Item(86, 2).Value = "Other Interactions"
$ws.Cells.Item(87, 2).Value = "Other Interactions"
$ws.Cells.Item(88, 2).Value = "Other Interactions"
$ws.Cells.Item(89, 2).Value = "Other Interactions"

# Fix Level 2 (column C) values: replace placeholder "Nan" with proper category, trim whitespace
$ws.Cells.Item(2, 3).Value = "Payments"
$ws.Cells.Item(3, 3).Value = "Payments"
$ws.Cells.Item(4, 3).Value = "Payments"
$ws.Cells.Item(5, 3).Value = "Payments"
$ws.Cells.Item(6, 3).Value = "Payments"
$ws.Cells.Item(7, 3).Value = "Payments"
$ws.Cells.Item(8, 3).Value = "Payments"
$ws.Cells.Item(9, 3).Value = "Payments"
$ws.Cells.Item(10, 3).Value = "Payments"
$ws.Cells.Item(11, 3).Value = "Payments"
$ws.Cells.Item(12, 3).Value = "Billing"
$ws.Cells.Item(13, 3).Value = "Pay Plans"
$ws.Cells.Item(14, 3).Value = "Account Balance"
$ws.Cells.Item(15, 3).Value = "Start"
$ws.Cells.Item(16, 3).Value = "Stop"
$ws.Cells.Item(17, 3).Value = "Transfer"
$ws.Cells.Item(18, 3).Value = "Outage"
$ws.Cells.Item(19, 3).Value = "Other Transactions"
$ws.Cells.Item(20, 3).Value = "Icm Technology Handled1"
$ws.Cells.Item(21, 3).Value = "Nonivr Technology2"
$ws.Cells.Item(22, 3).Value = "Handled Calls"
$ws.Cells.Item(23, 3).Value = "Transferred Calls4"
$ws.Cells.Item(24, 3).Value = "Abandoned Calls"
$ws.Cells.Item(25, 3).Value = "General"
$ws.Cells.Item(26, 3).Value = "Billing Balance"
$ws.Cells.Item(27, 3).Value = "Start"
$ws.Cells.Item(28, 3).Value = "Stop"
$ws.Cells.Item(29, 3).Value = "Transfer"
$ws.Cells.Item(30, 3).Value = "Bcsc (Business Customer Service Center)"
$ws.Cells.Item(31, 3).Value = "Emergency"
$ws.Cells.Item(32, 3).Value = "Outage"
$ws.Cells.Item(33, 3).Value = "Spanish (General)"
$ws.Cells.Item(34, 3).Value = "Payments"
$ws.Cells.Item(35, 3).Value = "Solar"
$ws.Cells.Item(36, 3).Value = "Specialty Lines"
$ws.Cells.Item(37, 3).Value = "Pay Plans"
$ws.Cells.Item(38, 3).Value = "Other Transactions"
$ws.Cells.Item(39, 3).Value = "Payment"
$ws.Cells.Item(42, 3).Value = "Billing"
$ws.Cells.Item(45, 3).Value = "Usage And Rates"
$ws.Cells.Item(46, 3).Value = "View Usage"
$ws.Cells.Item(51, 3).Value = "Outage"
$ws.Cells.Item(55, 3).Value = "Account Mgmt"
$ws.Cells.Item(56, 3).Value = "Change Billing Address"
$ws.Cells.Item(58, 3).Value = "Change User Name, Password"
$ws.Cells.Item(59, 3).Value = "Payment Account New"
$ws.Cells.Item(63, 3).Value = "Start"
$ws.Cells.Item(64, 3).Value = "Start"
$ws.Cells.Item(65, 3).Value = "Stop"
$ws.Cells.Item(66, 3).Value = "Transfer"
$ws.Cells.Item(67, 3).Value = "Pay Plans"
$ws.Cells.Item(70, 3).Value = "Care/Fera"
$ws.Cells.Item(71, 3).Value = "Alerts & Notifications"
$ws.Cells.Item(72, 3).Value = "Change Billing & Payment Alerts"
$ws.Cells.Item(81, 3).Value = "Pilot Light Appointments"
$ws.Cells.Item(85, 3).Value = "Budget Billing (Login)"
$ws.Cells.Item(86, 3).Value = "Other Interactions"
